$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "dumbi " to "Yuktha" in cell A7
$ws.Range("A7").Value = "Yuktha"

# Update the selection to match the post-edit state (single cell A7)
$ws.Range("A7").Select()
